# feat: add 2022-Q1 data
#
# 1. The sheet that used to be named "总计" (4th sheet, sheetId=4) becomes
#    "2022-Q1" and is populated with the new quarter's per-fund holdings
#    (it gains four extra columns: 基金规模/股票总仓位/仓位占比/仓位排名).
# 2. A brand-new "总计" summary sheet is appended right after it, rebuilt
#    from a copy of the original 4-column summary table so it inherits the
#    exact same sheet/cell formatting, with a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force the numeric-looking string to be stored as text (t="s") rather
    # than being auto-coerced into a number, then drop the "quote prefix"
    # text-format style that tags along, so the cell ends up unstyled -
    # matching its plain (un-bordered) neighbours.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Step 0: snapshot a pristine copy of the old "总计" sheet (still 4
#     columns wide) right after itself, BEFORE it gets renamed/widened, so
#     the new "总计" sheet naturally keeps the original A1:D layout/format
#     instead of inheriting the 2022-Q1 sheet's extra E:H columns. --------
$q1 = $wb.Worksheets.Item(4)
$q1.Name = "2022-Q1"
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item(5)
$total.Name = "总计"

# --- Step 1: turn "2022-Q1" into the new quarter's per-fund holdings -----
$q1.Cells.ClearContents()

# Header row: B1:D1 keep their original style; extend that same style to
# the new E1:H1 header cells before filling in their text.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 - 000593
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "000593"
Set-TextValue $q1.Range("C2") "易方达标普全球高端消费品指数增强(QDII)-美元现汇"
Set-TextValue $q1.Range("D2") "1.93"
Set-TextValue $q1.Range("E2") "92.46"
Set-TextValue $q1.Range("F2") "3.50"
Set-TextValue $q1.Range("G2") "0.0676"
$q1.Range("H2").Value = 10

# Row 3 - 005676
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "005676"
Set-TextValue $q1.Range("C3") "易方达标普全球高端消费品指数增强C(QDII) - 人民币"
Set-TextValue $q1.Range("D3") "1.93"
Set-TextValue $q1.Range("E3") "92.46"
Set-TextValue $q1.Range("F3") "3.50"
Set-TextValue $q1.Range("G3") "0.0676"
$q1.Range("H3").Value = 10

# Row 4 - 118002
$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "118002"
Set-TextValue $q1.Range("C4") "易方达标普全球高端消费品指数增强A(QDII) - 人民币"
Set-TextValue $q1.Range("D4") "1.93"
Set-TextValue $q1.Range("E4") "92.46"
Set-TextValue $q1.Range("F4") "3.50"
Set-TextValue $q1.Range("G4") "0.0676"
$q1.Range("H4").Value = 10

# --- Step 2: refresh the "总计" sheet with the rolling summary, now with
#     a new leading row for the just-added 2022-Q1 quarter. --------------
$total.Cells.ClearContents()

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.2

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.25

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.23

# Row 5 is new (the original summary table only ran to row 4), so its A5
# row-label cell needs the same style as A2:A4 carried over explicitly.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.23

# Restore the originally-active tab (first sheet) since none of the above
# navigation should leave the workbook's selection on a different sheet.
$wb.Worksheets.Item(1).Activate()
